$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.157353639602661
$ws.Range("B1").Value = 2.232425689697266
$ws.Range("C1").Value = 4.566696643829346
$ws.Range("D1").Value = 2.662504196166992
$ws.Range("E1").Value = 1.240180134773254
